$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- Column E (base): 18 cell(s) updated ---
$ws.Cells.Item(22, 5).Value = 'outputToCloud(resource)'
$ws.Cells.Item(23, 5).Value = 'prependText(var,prependWith)'
$ws.Cells.Item(24, 5).Value = 'repeatUntil(steps,maxWaitMs)'
$ws.Cells.Item(25, 5).Value = 'save(var,value)'
$ws.Cells.Item(26, 5).Value = 'saveCount(text,regex,saveVar)'
$ws.Cells.Item(27, 5).Value = 'saveMatches(text,regex,saveVar)'
$ws.Cells.Item(28, 5).Value = 'saveReplace(text,regex,replace,saveVar)'
$ws.Cells.Item(29, 5).Value = 'saveVariablesByPrefix(var,prefix)'
$ws.Cells.Item(30, 5).Value = 'saveVariablesByRegex(var,regex)'
$ws.Cells.Item(31, 5).Value = 'section(steps)'
$ws.Cells.Item(32, 5).Value = 'split(text,delim,saveVar)'
$ws.Cells.Item(33, 5).Value = 'startRecording()'
$ws.Cells.Item(34, 5).Value = 'stopRecording()'
$ws.Cells.Item(35, 5).Value = 'substringAfter(text,delim,saveVar)'
$ws.Cells.Item(36, 5).Value = 'substringBefore(text,delim,saveVar)'
$ws.Cells.Item(37, 5).Value = 'substringBetween(text,start,end,saveVar)'
$ws.Cells.Item(38, 5).Value = 'verbose(text)'
$ws.Cells.Item(39, 5).Value = 'waitFor(waitMs)'

# --- Column G (desktop): 94 cell(s) updated ---
$ws.Cells.Item(5, 7).Value = 'assertElementNotPresent(name)'
$ws.Cells.Item(6, 7).Value = 'assertElementPresent(name)'
$ws.Cells.Item(7, 7).Value = 'assertEnabled(name)'
$ws.Cells.Item(8, 7).Value = 'assertHierCells(matchBy,column,expected,nestedOnly)'
$ws.Cells.Item(9, 7).Value = 'assertHierRow(matchBy,expected)'
$ws.Cells.Item(10, 7).Value = 'assertListCount(count)'
$ws.Cells.Item(11, 7).Value = 'assertLocatorNotPresent(locator)'
$ws.Cells.Item(12, 7).Value = 'assertLocatorPresent(locator)'
$ws.Cells.Item(13, 7).Value = 'assertMenuEnabled(menu)'
$ws.Cells.Item(14, 7).Value = 'assertModalDialogNotPresent()'
$ws.Cells.Item(15, 7).Value = 'assertModalDialogPresent()'
$ws.Cells.Item(16, 7).Value = 'assertModalDialogTitle(title)'
$ws.Cells.Item(17, 7).Value = 'assertModalDialogTitleByLocator(locator,title)'
$ws.Cells.Item(18, 7).Value = 'assertNotChecked(name)'
$ws.Cells.Item(19, 7).Value = 'assertSelected(name,text)'
$ws.Cells.Item(20, 7).Value = 'assertTableCell(row,column,contains)'
$ws.Cells.Item(21, 7).Value = 'assertTableColumnContains(column,contains)'
$ws.Cells.Item(22, 7).Value = 'assertTableContains(contains)'
$ws.Cells.Item(23, 7).Value = 'assertTableRowContains(row,contains)'
$ws.Cells.Item(24, 7).Value = 'assertText(name,expected)'
$ws.Cells.Item(25, 7).Value = 'assertWindowTitleContains(contains)'
$ws.Cells.Item(26, 7).Value = 'clear(locator)'
$ws.Cells.Item(27, 7).Value = 'clearCombo(name)'
$ws.Cells.Item(28, 7).Value = 'clearModalDialog(var,button)'
$ws.Cells.Item(29, 7).Value = 'clearTextArea(name)'
$ws.Cells.Item(30, 7).Value = 'clearTextBox(name)'
$ws.Cells.Item(31, 7).Value = 'clickButton(name)'
$ws.Cells.Item(32, 7).Value = 'clickByLocator(locator)'
$ws.Cells.Item(33, 7).Value = 'clickCheckBox(name)'
$ws.Cells.Item(34, 7).Value = 'clickElementOffset(name,xOffset,yOffset)'
$ws.Cells.Item(35, 7).Value = 'clickExplorerBar(group,item)'
$ws.Cells.Item(36, 7).Value = 'clickFirstMatchRow(nameValues)'
$ws.Cells.Item(37, 7).Value = 'clickFirstMatchedList(contains)'
$ws.Cells.Item(38, 7).Value = 'clickIcon(label)'
$ws.Cells.Item(39, 7).Value = 'clickList(row)'
$ws.Cells.Item(40, 7).Value = 'clickMenu(menu)'
$ws.Cells.Item(41, 7).Value = 'clickOffset(locator,xOffset,yOffset)'
$ws.Cells.Item(42, 7).Value = 'clickRadio(name)'
$ws.Cells.Item(43, 7).Value = 'clickScreen(button,modifiers,x,y)'
$ws.Cells.Item(44, 7).Value = 'clickTab(group,name)'
$ws.Cells.Item(45, 7).Value = 'clickTableCell(row,column)'
$ws.Cells.Item(46, 7).Value = 'clickTableRow(row)'
$ws.Cells.Item(47, 7).Value = 'clickTextPane(name,criteria)'
$ws.Cells.Item(48, 7).Value = 'clickTextPaneRow(var,index)'
$ws.Cells.Item(49, 7).Value = 'closeApplication()'
$ws.Cells.Item(50, 7).Value = 'collapseHierTable()'
$ws.Cells.Item(51, 7).Value = 'editCurrentRow(nameValues)'
$ws.Cells.Item(52, 7).Value = 'editHierCells(var,matchBy,nameValues)'
$ws.Cells.Item(53, 7).Value = 'editTableCells(row,nameValues)'
$ws.Cells.Item(54, 7).Value = 'getRowCount(var)'
$ws.Cells.Item(55, 7).Value = 'hideExplorerBar()'
$ws.Cells.Item(56, 7).Value = 'login(form,username,password)'
$ws.Cells.Item(57, 7).Value = 'maximize()'
$ws.Cells.Item(58, 7).Value = 'minimize()'
$ws.Cells.Item(59, 7).Value = 'mouseWheel(amount,modifiers,x,y)'
$ws.Cells.Item(60, 7).Value = 'resize(width,height)'
$ws.Cells.Item(61, 7).Value = 'saveAllTableRows(var)'
$ws.Cells.Item(62, 7).Value = 'saveAttributeByLocator(var,locator,attribute)'
$ws.Cells.Item(63, 7).Value = 'saveElementCount(var,name)'
$ws.Cells.Item(64, 7).Value = 'saveFirstListData(var,contains)'
$ws.Cells.Item(65, 7).Value = 'saveFirstMatchedListIndex(var,contains)'
$ws.Cells.Item(66, 7).Value = 'saveHierCells(var,matchBy,column,nestedOnly)'
$ws.Cells.Item(67, 7).Value = 'saveHierRow(var,matchBy)'
$ws.Cells.Item(68, 7).Value = 'saveListData(var,contains)'
$ws.Cells.Item(69, 7).Value = 'saveLocatorCount(var,locator)'
$ws.Cells.Item(70, 7).Value = 'saveModalDialogText(var)'
$ws.Cells.Item(71, 7).Value = 'saveModalDialogTextByLocator(var,locator)'
$ws.Cells.Item(72, 7).Value = 'saveProcessId(var,locator)'
$ws.Cells.Item(73, 7).Value = 'saveRowCount(var)'
$ws.Cells.Item(74, 7).Value = 'saveTableRows(var,contains)'
$ws.Cells.Item(75, 7).Value = 'saveTableRowsRange(var,beginRow,endRow)'
$ws.Cells.Item(76, 7).Value = 'saveText(var,name)'
$ws.Cells.Item(77, 7).Value = 'saveTextByLocator(var,locator)'
$ws.Cells.Item(78, 7).Value = 'saveTextPane(var,name,criteria)'
$ws.Cells.Item(79, 7).Value = 'saveWindowTitle(var)'
$ws.Cells.Item(80, 7).Value = 'scanTable(var,name)'
$ws.Cells.Item(81, 7).Value = 'selectCombo(name,text)'
$ws.Cells.Item(82, 7).Value = 'sendKeysToTextBox(name,text1,text2,text3,text4)'
$ws.Cells.Item(83, 7).Value = 'showExplorerBar()'
$ws.Cells.Item(84, 7).Value = 'toggleExplorerBar()'
$ws.Cells.Item(85, 7).Value = 'typeAppendTextArea(name,text1,text2,text3,text4)'
$ws.Cells.Item(86, 7).Value = 'typeAppendTextBox(name,text1,text2,text3,text4)'
$ws.Cells.Item(87, 7).Value = 'typeByLocator(locator,text)'
$ws.Cells.Item(88, 7).Value = 'typeKeys(os,keystrokes)'
$ws.Cells.Item(89, 7).Value = 'typeTextArea(name,text1,text2,text3,text4)'
$ws.Cells.Item(90, 7).Value = 'typeTextBox(name,text1,text2,text3,text4)'
$ws.Cells.Item(91, 7).Value = 'useApp(appId)'
$ws.Cells.Item(92, 7).Value = 'useForm(formName)'
$ws.Cells.Item(93, 7).Value = 'useHierTable(var,name)'
$ws.Cells.Item(94, 7).Value = 'useList(var,name)'
$ws.Cells.Item(95, 7).Value = 'useTable(var,name)'
$ws.Cells.Item(96, 7).Value = 'useTableRow(var,row)'
$ws.Cells.Item(97, 7).Value = 'waitFor(name,maxWaitMs)'
$ws.Cells.Item(98, 7).Value = 'waitForLocator(locator,maxWaitMs)'

# --- Column I (external): 1 cell(s) updated ---
$ws.Cells.Item(5, 9).Value = 'tail(id,file)'

# --- Column J (image): 1 cell(s) updated ---
$ws.Cells.Item(7, 10).Value = 'saveDiff(var,baseline,actual)'

# --- Column M (json): 7 cell(s) updated ---
$ws.Cells.Item(12, 13).Value = 'compact(var,json,removeEmpty)'
$ws.Cells.Item(13, 13).Value = 'fromCsv(csv,header,jsonFile)'
$ws.Cells.Item(14, 13).Value = 'minify(json,var)'
$ws.Cells.Item(15, 13).Value = 'storeCount(json,jsonpath,var)'
$ws.Cells.Item(16, 13).Value = 'storeKeys(json,jsonpath,var)'
$ws.Cells.Item(17, 13).Value = 'storeValue(json,jsonpath,var)'
$ws.Cells.Item(18, 13).Value = 'storeValues(json,jsonpath,var)'

# --- Column Y (web): 111 cell(s) updated ---
$ws.Cells.Item(24, 25).Value = 'assertMultiSelect(locator)'
$ws.Cells.Item(25, 25).Value = 'assertNotChecked(locator)'
$ws.Cells.Item(26, 25).Value = 'assertNotFocus(locator)'
$ws.Cells.Item(27, 25).Value = 'assertNotText(locator,text)'
$ws.Cells.Item(28, 25).Value = 'assertNotVisible(locator)'
$ws.Cells.Item(29, 25).Value = 'assertOneMatch(locator)'
$ws.Cells.Item(30, 25).Value = 'assertScrollbarHNotPresent(locator)'
$ws.Cells.Item(31, 25).Value = 'assertScrollbarHPresent(locator)'
$ws.Cells.Item(32, 25).Value = 'assertScrollbarVNotPresent(locator)'
$ws.Cells.Item(33, 25).Value = 'assertScrollbarVPresent(locator)'
$ws.Cells.Item(34, 25).Value = 'assertSingleSelect(locator)'
$ws.Cells.Item(35, 25).Value = 'assertTable(locator,row,column,text)'
$ws.Cells.Item(36, 25).Value = 'assertText(locator,text)'
$ws.Cells.Item(37, 25).Value = 'assertTextContains(locator,text)'
$ws.Cells.Item(38, 25).Value = 'assertTextCount(locator,text,count)'
$ws.Cells.Item(39, 25).Value = 'assertTextList(locator,list,ignoreOrder)'
$ws.Cells.Item(40, 25).Value = 'assertTextMatches(text,minMatch,scrollTo)'
$ws.Cells.Item(41, 25).Value = 'assertTextNotContains(locator,text)'
$ws.Cells.Item(42, 25).Value = 'assertTextNotPresent(text)'
$ws.Cells.Item(43, 25).Value = 'assertTextOrder(locator,descending)'
$ws.Cells.Item(44, 25).Value = 'assertTextPresent(text)'
$ws.Cells.Item(45, 25).Value = 'assertTitle(text)'
$ws.Cells.Item(46, 25).Value = 'assertValue(locator,value)'
$ws.Cells.Item(47, 25).Value = 'assertValueOrder(locator,descending)'
$ws.Cells.Item(48, 25).Value = 'assertVisible(locator)'
$ws.Cells.Item(49, 25).Value = 'checkAll(locator)'
$ws.Cells.Item(50, 25).Value = 'clearLocalStorage()'
$ws.Cells.Item(51, 25).Value = 'click(locator)'
$ws.Cells.Item(52, 25).Value = 'clickAll(locator)'
$ws.Cells.Item(53, 25).Value = 'clickAndWait(locator,waitMs)'
$ws.Cells.Item(54, 25).Value = 'clickByLabel(label)'
$ws.Cells.Item(55, 25).Value = 'clickByLabelAndWait(label,waitMs)'
$ws.Cells.Item(56, 25).Value = 'clickOffset(locator,x,y)'
$ws.Cells.Item(57, 25).Value = 'clickWithKeys(locator,keys)'
$ws.Cells.Item(58, 25).Value = 'close()'
$ws.Cells.Item(59, 25).Value = 'closeAll()'
$ws.Cells.Item(60, 25).Value = 'deselect(locator,text)'
$ws.Cells.Item(61, 25).Value = 'deselectMulti(locator,array)'
$ws.Cells.Item(62, 25).Value = 'dismissInvalidCert()'
$ws.Cells.Item(63, 25).Value = 'dismissInvalidCertPopup()'
$ws.Cells.Item(64, 25).Value = 'doubleClick(locator)'
$ws.Cells.Item(65, 25).Value = 'doubleClickAndWait(locator,waitMs)'
$ws.Cells.Item(66, 25).Value = 'doubleClickByLabel(label)'
$ws.Cells.Item(67, 25).Value = 'doubleClickByLabelAndWait(label,waitMs)'
$ws.Cells.Item(68, 25).Value = 'dragAndDrop(fromLocator,toLocator)'
$ws.Cells.Item(69, 25).Value = 'dragTo(fromLocator,xOffset,yOffset)'
$ws.Cells.Item(70, 25).Value = 'editLocalStorage(key,value)'
$ws.Cells.Item(71, 25).Value = 'executeScript(var,script)'
$ws.Cells.Item(72, 25).Value = 'focus(locator)'
$ws.Cells.Item(73, 25).Value = 'goBack()'
$ws.Cells.Item(74, 25).Value = 'goBackAndWait()'
$ws.Cells.Item(75, 25).Value = 'maximizeWindow()'
$ws.Cells.Item(76, 25).Value = 'mouseOver(locator)'
$ws.Cells.Item(77, 25).Value = 'open(url)'
$ws.Cells.Item(78, 25).Value = 'openAndWait(url,waitMs)'
$ws.Cells.Item(79, 25).Value = 'openHttpBasic(url,username,password)'
$ws.Cells.Item(80, 25).Value = 'openIgnoreTimeout(url)'
$ws.Cells.Item(81, 25).Value = 'refresh()'
$ws.Cells.Item(82, 25).Value = 'refreshAndWait()'
$ws.Cells.Item(83, 25).Value = 'resizeWindow(width,height)'
$ws.Cells.Item(84, 25).Value = 'rightClick(locator)'
$ws.Cells.Item(85, 25).Value = 'saveAllWindowIds(var)'
$ws.Cells.Item(86, 25).Value = 'saveAllWindowNames(var)'
$ws.Cells.Item(87, 25).Value = 'saveAttribute(var,locator,attrName)'
$ws.Cells.Item(88, 25).Value = 'saveAttributeList(var,locator,attrName)'
$ws.Cells.Item(89, 25).Value = 'saveCount(var,locator)'
$ws.Cells.Item(90, 25).Value = 'saveDivsAsCsv(headers,rows,cells,nextPage,file)'
$ws.Cells.Item(91, 25).Value = 'saveElement(var,locator)'
$ws.Cells.Item(92, 25).Value = 'saveElements(var,locator)'
$ws.Cells.Item(93, 25).Value = 'saveInfiniteDivsAsCsv(config,file)'
$ws.Cells.Item(94, 25).Value = 'saveInfiniteTableAsCsv(config,file)'
$ws.Cells.Item(95, 25).Value = 'saveLocalStorage(var,key)'
$ws.Cells.Item(96, 25).Value = 'saveLocation(var)'
$ws.Cells.Item(97, 25).Value = 'savePageAs(var,sessionIdName,url)'
$ws.Cells.Item(98, 25).Value = 'savePageAsFile(sessionIdName,url,file)'
$ws.Cells.Item(99, 25).Value = 'saveTableAsCsv(locator,nextPageLocator,file)'
$ws.Cells.Item(100, 25).Value = 'saveText(var,locator)'
$ws.Cells.Item(101, 25).Value = 'saveTextArray(var,locator)'
$ws.Cells.Item(102, 25).Value = 'saveTextSubstringAfter(var,locator,delim)'
$ws.Cells.Item(103, 25).Value = 'saveTextSubstringBefore(var,locator,delim)'
$ws.Cells.Item(104, 25).Value = 'saveTextSubstringBetween(var,locator,start,end)'
$ws.Cells.Item(105, 25).Value = 'saveValue(var,locator)'
$ws.Cells.Item(106, 25).Value = 'saveValues(var,locator)'
$ws.Cells.Item(107, 25).Value = 'scrollElement(locator,xOffset,yOffset)'
$ws.Cells.Item(108, 25).Value = 'scrollLeft(locator,pixel)'
$ws.Cells.Item(109, 25).Value = 'scrollPage(xOffset,yOffset)'
$ws.Cells.Item(110, 25).Value = 'scrollRight(locator,pixel)'
$ws.Cells.Item(111, 25).Value = 'scrollTo(locator)'
$ws.Cells.Item(112, 25).Value = 'select(locator,text)'
$ws.Cells.Item(113, 25).Value = 'selectFrame(locator)'
$ws.Cells.Item(114, 25).Value = 'selectMulti(locator,array)'
$ws.Cells.Item(115, 25).Value = 'selectMultiOptions(locator)'
$ws.Cells.Item(116, 25).Value = 'selectText(locator)'
$ws.Cells.Item(117, 25).Value = 'selectWindow(winId)'
$ws.Cells.Item(118, 25).Value = 'selectWindowAndWait(winId,waitMs)'
$ws.Cells.Item(119, 25).Value = 'selectWindowByIndex(index)'
$ws.Cells.Item(120, 25).Value = 'selectWindowByIndexAndWait(index,waitMs)'
$ws.Cells.Item(121, 25).Value = 'toggleSelections(locator)'
$ws.Cells.Item(122, 25).Value = 'type(locator,value)'
$ws.Cells.Item(123, 25).Value = 'typeKeys(locator,value)'
$ws.Cells.Item(124, 25).Value = 'uncheckAll(locator)'
$ws.Cells.Item(125, 25).Value = 'unselectAllText()'
$ws.Cells.Item(126, 25).Value = 'updateAttribute(locator,attrName,value)'
$ws.Cells.Item(127, 25).Value = 'upload(fieldLocator,file)'
$ws.Cells.Item(128, 25).Value = 'verifyContainText(locator,text)'
$ws.Cells.Item(129, 25).Value = 'verifyText(locator,text)'
$ws.Cells.Item(130, 25).Value = 'wait(waitMs)'
$ws.Cells.Item(131, 25).Value = 'waitForElementPresent(locator)'
$ws.Cells.Item(132, 25).Value = 'waitForPopUp(winId,waitMs)'
$ws.Cells.Item(133, 25).Value = 'waitForTextPresent(text)'
$ws.Cells.Item(134, 25).Value = 'waitForTitle(text)'

# --- Column AD (xml): 21 cell(s) updated ---
$ws.Cells.Item(7, 30).Value = 'assertSoap(wsdl,xml)'
$ws.Cells.Item(8, 30).Value = 'assertSoapFaultCode(expected,xml)'
$ws.Cells.Item(9, 30).Value = 'assertSoapFaultString(expected,xml)'
$ws.Cells.Item(10, 30).Value = 'assertValue(xml,xpath,expected)'
$ws.Cells.Item(11, 30).Value = 'assertValues(xml,xpath,array,exactOrder)'
$ws.Cells.Item(12, 30).Value = 'assertWellformed(xml)'
$ws.Cells.Item(13, 30).Value = 'beautify(xml,var)'
$ws.Cells.Item(14, 30).Value = 'clear(xml,xpath,var)'
$ws.Cells.Item(15, 30).Value = 'delete(xml,xpath,var)'
$ws.Cells.Item(16, 30).Value = 'insertAfter(xml,xpath,content,var)'
$ws.Cells.Item(17, 30).Value = 'insertBefore(xml,xpath,content,var)'
$ws.Cells.Item(18, 30).Value = 'minify(xml,var)'
$ws.Cells.Item(19, 30).Value = 'prepend(xml,xpath,content,var)'
$ws.Cells.Item(20, 30).Value = 'replace(xml,xpath,content,var)'
$ws.Cells.Item(21, 30).Value = 'replaceIn(xml,xpath,content,var)'
$ws.Cells.Item(22, 30).Value = 'storeCount(xml,xpath,var)'
$ws.Cells.Item(23, 30).Value = 'storeSoapFaultCode(var,xml)'
$ws.Cells.Item(24, 30).Value = 'storeSoapFaultDetail(var,xml)'
$ws.Cells.Item(25, 30).Value = 'storeSoapFaultString(var,xml)'
$ws.Cells.Item(26, 30).Value = 'storeValue(xml,xpath,var)'
$ws.Cells.Item(27, 30).Value = 'storeValues(xml,xpath,var)'

# --- Update defined name ranges to reflect new list extents ---
$wb.Names.Item('base').RefersTo = "'#system'!\$E\$2:\$E\$39"
$wb.Names.Item('desktop').RefersTo = "'#system'!\$G\$2:\$G\$98"
$wb.Names.Item('external').RefersTo = "'#system'!\$I\$2:\$I\$5"
$wb.Names.Item('image').RefersTo = "'#system'!\$J\$2:\$J\$7"
$wb.Names.Item('json').RefersTo = "'#system'!\$M\$2:\$M\$18"
$wb.Names.Item('web').RefersTo = "'#system'!\$Y\$2:\$Y\$134"
$wb.Names.Item('xml').RefersTo = "'#system'!\$AD\$2:\$AD\$27"
